# RF classify (read data header bug)
#
# The sheet used to keep a second, duplicate copy of the raw "read data"
# header/columns (H:M) next to the computed phase table (A:F), plus a
# one-off "Video start time:" row (row 9) that the A/B formulas referenced
# via the fixed anchor B$9. That anchor was wrong (a read-data-header bug),
# so this edit:
#   1) freezes the (already-correct, cached) computed start/end times in
#      A2:B7 to plain values - they no longer depend on the raw H:I/B9
#      helper cells,
#   2) removes the duplicate raw-data header block H1:M7,
#   3) removes the now-unused "Video start time:" row (old row 9), which
#      shifts the trailing helper cell from B11 up to B10,
#   4) leaves the selection where the author left it (R13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Freeze the computed start/end-time formulas in A2:B7 to their current
#    values so they no longer depend on the raw H:I columns / B9 anchor.
for ($r = 2; $r -le 7; $r++) {
    $aVal = $ws.Cells.Item($r, 1).Value2
    $bVal = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value = $aVal
    $ws.Cells.Item($r, 2).Value = $bVal
}

# 2) Remove the duplicate raw "read data" header/table in columns H:M.
$ws.Range("H1:M7").Delete()

# 3) Remove the stray "Video start time:" row (old row 9). This shifts the
#    trailing helper value up from B11 to B10.
$ws.Rows("9").Delete()

# 4) Restore the author's last selection.
$ws.Range("R13").Select()
